$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 206
$ws.Cells.Item(2, 2).Value = 14
$ws.Cells.Item(2, 3).Value = "許*綸"
$ws.Cells.Item(2, 4).Value = "2024-03-02 06:48:46"
$ws.Cells.Item(2, 5).Value = "2024-03-02 06:49:23"
$ws.Cells.Item(2, 6).Value = "OUT"
